$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) MD410 Attendance: insert a new registree row at row 222
#    ("tyler", "natascha", "Kensington", "No", "No", "410E"), pushing the
#    existing rows down by one.
# ---------------------------------------------------------------------------
$wsMD410 = $wb.Worksheets.Item("MD410 Attendance")

$wsMD410.Rows.Item(222).Insert()
$wsMD410.Range("A223:F223").Copy()
$wsMD410.Range("A222:F222").PasteSpecial(-4122)
$wsMD410.Rows.Item(222).RowHeight = 25

$wsMD410.Range("A222").Value = "tyler"
$wsMD410.Range("B222").Value = "natascha"
$wsMD410.Range("C222").Value = "Kensington"
$wsMD410.Range("D222").Value = "No"
$wsMD410.Range("E222").Value = "No"
$wsMD410.Range("F222").Value = "410E"

# Bump the "Number of attendees" summary line (now shifted to row 241).
$wsMD410.Range("A241").Value = "Number of attendees: 238"

# Update the report-generation timestamp in the title cell.
$wsMD410.Range("A1").Value = "MD410 Registrees as of 22/04/2021 13:11"

# ---------------------------------------------------------------------------
# 2) 410E Attendance: same new registree, inserted at row 118 (this sheet
#    has no District column, so only columns A:E).
# ---------------------------------------------------------------------------
$ws410E = $wb.Worksheets.Item("410E Attendance")

$ws410E.Rows.Item(118).Insert()
$ws410E.Range("A119:E119").Copy()
$ws410E.Range("A118:E118").PasteSpecial(-4122)
$ws410E.Rows.Item(118).RowHeight = 25

$ws410E.Range("A118").Value = "tyler"
$ws410E.Range("B118").Value = "natascha"
$ws410E.Range("C118").Value = "Kensington"
$ws410E.Range("D118").Value = "No"
$ws410E.Range("E118").Value = "No"

# Bump the "Number of attendees" summary line (now shifted to row 130).
$ws410E.Range("A130").Value = "Number of attendees: 127"

$ws410E.Range("A1").Value = "410E Registrees as of 22/04/2021 13:11"

# ---------------------------------------------------------------------------
# 3) 410W Attendance: no roster change (the new registree is 410E), just the
#    refreshed report timestamp.
# ---------------------------------------------------------------------------
$ws410W = $wb.Worksheets.Item("410W Attendance")
$ws410W.Range("A1").Value = "410W Registrees as of 22/04/2021 13:11"

# ---------------------------------------------------------------------------
# 4) 410E Voting / 410W Voting: timestamp-only refresh (the new registree is
#    not a voter, so the voting tallies are unaffected).
# ---------------------------------------------------------------------------
$ws410EVoting = $wb.Worksheets.Item("410E Voting")
$ws410EVoting.Range("A1").Value = "410E Voting details as of 22/04/2021 13:11"

$ws410WVoting = $wb.Worksheets.Item("410W Voting")
$ws410WVoting.Range("A1").Value = "410W Voting details as of 22/04/2021 13:11"
